$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-5 (the original rows 6-9 are removed entirely, and the
# remaining rows 2-5 get updated values reflecting the refreshed TPM numbers).
$data = @(
    @{ A="ECs";           D="Resolving-Mac"; G=4.931298666666667; H=14.793896;       I=0.243258826942244;  J=0.243258826942244;  K=3; L=1; M=15.40797466666667; N=46.223924; O=1; P=1; Q=75.98132492976711;  R=683.831924367904;  S=0.243258826942244;  T=0.243258826942244 },
    @{ A="FAPs";          D="Resolving-Mac"; G=2.442036;          H=7.326108;         I=0.1204645779673042; J=0.1204645779673042; K=3; L=1; M=15.40797466666667; N=46.223924; O=1; P=1; Q=37.626828823088;    R=338.641459407792;  S=0.1204645779673042; T=0.1204645779673042 },
    @{ A="MuSCs";         D="Resolving-Mac"; G=3.440487666666666; H=10.321463;        I=0.1697177661454274; J=0.1697177661454274; K=3; L=1; M=15.40797466666667; N=46.223924; O=1; P=1; Q=53.01094680897911;  R=477.098521280812;  S=0.1697177661454274; T=0.1697177661454274 },
    @{ A="Resolving-Mac"; D="Resolving-Mac"; G=9.457995666666667; H=28.373987;        I=0.4665588289450244; J=0.4665588289450244; K=3; L=1; M=15.40797466666667; N=46.223924; O=1; P=1; Q=145.7285576294431;  R=1311.557018664988; S=0.4665588289450244; T=0.4665588289450244 }
)

$rowIdx = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIdx, 1).Value = $rec.A        # A: Sending cluster
    $ws.Cells.Item($rowIdx, 2).Value = "St6gal1"      # B: Ligand symbol
    $ws.Cells.Item($rowIdx, 3).Value = "Cd22"         # C: Receptor symbol
    $ws.Cells.Item($rowIdx, 4).Value = $rec.D         # D: Target cluster
    $ws.Cells.Item($rowIdx, 5).Value = 3               # E: Ligand-expressing cells
    $ws.Cells.Item($rowIdx, 6).Value = 1               # F: Ligand detection rate
    $ws.Cells.Item($rowIdx, 7).Value = $rec.G
    $ws.Cells.Item($rowIdx, 8).Value = $rec.H
    $ws.Cells.Item($rowIdx, 9).Value = $rec.I
    $ws.Cells.Item($rowIdx, 10).Value = $rec.J
    $ws.Cells.Item($rowIdx, 11).Value = $rec.K
    $ws.Cells.Item($rowIdx, 12).Value = $rec.L
    $ws.Cells.Item($rowIdx, 13).Value = $rec.M
    $ws.Cells.Item($rowIdx, 14).Value = $rec.N
    $ws.Cells.Item($rowIdx, 15).Value = $rec.O
    $ws.Cells.Item($rowIdx, 16).Value = $rec.P
    $ws.Cells.Item($rowIdx, 17).Value = $rec.Q
    $ws.Cells.Item($rowIdx, 18).Value = $rec.R
    $ws.Cells.Item($rowIdx, 19).Value = $rec.S
    $ws.Cells.Item($rowIdx, 20).Value = $rec.T
    $rowIdx++
}

# Remove the old trailing rows (6-9) that no longer exist in the refreshed data.
$ws.Range("A6:T9").EntireRow.Delete() | Out-Null
